$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Append 4 new sales rows (10-13) below the existing table (rows 1-9).
# New literal text values ("2023-10-31", "2400", "5400") must land in
# the sheet as genuine shared strings (t="s"), not auto-converted
# numbers/dates, and the pre-existing repeated values ("121210",
# "Oral B", "800", "121212", "Caro White Cream", "1800", "3") must
# reuse the workbook's existing shared-string entries rather than
# creating duplicates.
#
# Technique: build each brand-new text value once in a scratch cell
# that has been explicitly formatted as Text ("@"), copy it, and
# paste-special "values only" into the destination - this carries
# over the string flavour without carrying the scratch cell's number
# format into the destination. The scratch cells are fully cleared
# afterwards so no stray content/format is left behind.
# ------------------------------------------------------------------

$scratchDate = $ws.Range("H1")
$scratchDate.NumberFormat = "@"
$scratchDate.Value = "2023-10-31"

$scratchAmt1 = $ws.Range("H2")
$scratchAmt1.NumberFormat = "@"
$scratchAmt1.Value = "2400"

$scratchAmt2 = $ws.Range("H3")
$scratchAmt2.NumberFormat = "@"
$scratchAmt2.Value = "5400"

# --- Row 10: duplicate of row 2 (PC 121210 / Oral B) with new date,
#     qty 3, amt 2400 -------------------------------------------------
$ws.Range("A2:F2").Copy()
$ws.Range("A10:F10").PasteSpecial(-4163)

$scratchDate.Copy()
$ws.Range("B10").PasteSpecial(-4163)

$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4163)

$scratchAmt1.Copy()
$ws.Range("F10").PasteSpecial(-4163)

# --- Rows 11-13: three duplicates of row 3 (PC 121212 / Caro White
#     Cream) with new date, qty 3, amt 5400 ---------------------------
foreach ($r in 11, 12, 13) {
    $ws.Range("A3:F3").Copy()
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4163)

    $scratchDate.Copy()
    $ws.Range("B" + $r).PasteSpecial(-4163)

    $ws.Range("E9").Copy()
    $ws.Range("E" + $r).PasteSpecial(-4163)

    $scratchAmt2.Copy()
    $ws.Range("F" + $r).PasteSpecial(-4163)
}

# --- Clean up scratch cells so nothing extra is left in the sheet ---
$ws.Range("H1:H3").Clear()
$excel.CutCopyMode = $false
